# Apply resource-assignment updates to the "Tabla de actividades" sheet
# (commit: "modificacion de gantt y asignacion de recursos en gantt")
#
# Column F ("Recursos ") gets additional resource codes (S[xx] / T[xx])
# appended to the already-assigned P[xx] (programmer) codes for a number
# of activities.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = "P[01],S[02],T[02]"                 # ADA01001
$ws.Range("F8").Value  = "P[04],P[03],S[02],T[03],T[01]"     # ADA01006
$ws.Range("F9").Value  = "P[01],P[02],S[01],T[02],T[04]"     # ADA01007
$ws.Range("F26").Value = "P[03],S[01],T[03]"                 # BD01002
$ws.Range("F27").Value = "P[03],S[01],T[03]"                 # BD01003
$ws.Range("F28").Value = "P[02],S[02],T[04]"                 # BD01004
$ws.Range("F29").Value = "P[04],S[02],T[01]"                 # BD01005
$ws.Range("F42").Value = " "                                 # Prog01001
$ws.Range("F51").Value = "P[03],S[01],T[03],"                # Proy01001
$ws.Range("F53").Value = "P[01],S[02],T[02]"                 # Proy01003
$ws.Range("F76").Value = "P[03],S[11],T[03]"                 # SO01006

# Update the sheet view to reflect where the author was working (zoomed
# in further down the Gantt, around row 25 / column F)
$excel.ActiveWindow.Zoom = 62
$ws.Range("F80").Select()
